$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a bare empty paragraph ("<w:p/>") right after the paragraph
#    that reads "In jetpack, we use paddings for margin" (and before the
#    following empty Heading1 paragraph).
# ---------------------------------------------------------------------
$findRng = $d.Content
$found = $findRng.Find.Execute("In jetpack, we use paddings for margin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find paragraph 'In jetpack, we use paddings for margin'"
}
$jetpackIndex = $findRng.Paragraphs.Item(1).Index
$jetpackPara = $d.Paragraphs.Item($jetpackIndex)

# Use a plain Range positioned exactly at the end of the paragraph (not a
# Collapse()'d duplicate) and InsertXML a genuinely empty <w:p/> so no
# stray empty run is produced and the preceding paragraph's own text is
# left untouched.
$afterJetpack = $d.Range($jetpackPara.Range.End, $jetpackPara.Range.End)
[void]$afterJetpack.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

# ---------------------------------------------------------------------
# 2) Insert four new plain paragraphs right after the
#    "Modifier Attributes on Text" heading paragraph.
# ---------------------------------------------------------------------
$findRng2 = $d.Content
$found2 = $findRng2.Find.Execute("Modifier Attributes on Text", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find paragraph 'Modifier Attributes on Text'"
}
$headingIndex = $findRng2.Paragraphs.Item(1).Index

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$newParaTexts = @(
    "Text($openQuote" + "Talha$closeQuote, modifier = Modifier.offset(50.dp, 20.dp)",
    "50.dp -> right margin",
    "20.dp -> bottom",
    "It does not push other elements, it only takes given elements."
)

$curIndex = $headingIndex
foreach ($txt in $newParaTexts) {
    $curPara = $d.Paragraphs.Item($curIndex)
    $insertPoint = $curPara.Range
    $insertPoint.Collapse(0)
    $insertPoint.InsertParagraphAfter()
    $curIndex = $curIndex + 1
    $newPara = $d.Paragraphs.Item($curIndex)
    $newPara.Style = "Normal"
    $newPara.Range.Text = $txt
}
